$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0.4480947554111481
$ws.Cells.Item(3, 4).Value = 0.3647661209106445
$ws.Cells.Item(4, 4).Value = 0.8815181255340576
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 0.4906995892524719
$ws.Cells.Item(6, 4).Value = 0.03579447790980339
$ws.Cells.Item(7, 4).Value = 0.8550922870635986
$ws.Cells.Item(8, 4).Value = 0.6009546518325806
$ws.Cells.Item(9, 4).Value = 0.2650291323661804
$ws.Cells.Item(10, 4).Value = 0.6092824935913086
$ws.Cells.Item(11, 4).Value = 0.2607499063014984
$ws.Cells.Item(12, 4).Value = 0.3480855822563171
$ws.Cells.Item(13, 4).Value = 0.9318639039993286
$ws.Cells.Item(14, 4).Value = 0.8949315547943115
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0.4818087816238403
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 0.2817938923835754
$ws.Cells.Item(17, 4).Value = 0.6917387843132019
$ws.Cells.Item(18, 4).Value = 0.6596561670303345
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 0.4472375214099884
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 0.3857227563858032
$ws.Cells.Item(21, 4).Value = 0.3868101835250854
$ws.Cells.Item(22, 4).Value = 0.6083816289901733
$ws.Cells.Item(23, 4).Value = 0.8168272376060486
$ws.Cells.Item(24, 4).Value = 0.4704257249832153
$ws.Cells.Item(25, 4).Value = 0.6048879623413086
$ws.Cells.Item(26, 4).Value = 0.8973966836929321
$ws.Cells.Item(27, 4).Value = 0.09893088042736053
$ws.Cells.Item(28, 4).Value = 0.7424843311309814
$ws.Cells.Item(29, 4).Value = 0.2871742844581604
$ws.Cells.Item(30, 4).Value = 0.7195528745651245
$ws.Cells.Item(31, 4).Value = 0.6281295418739319
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 0.4695977568626404
$ws.Cells.Item(33, 4).Value = 0.7509008049964905
$ws.Cells.Item(34, 4).Value = 0.8773109316825867
$ws.Cells.Item(35, 4).Value = 0.6579916477203369
$ws.Cells.Item(36, 4).Value = 0.1931407898664474
$ws.Cells.Item(37, 4).Value = 0.6744139194488525
$ws.Cells.Item(38, 4).Value = 0.6589869260787964
$ws.Cells.Item(39, 4).Value = 0.9035844206809998
$ws.Cells.Item(40, 4).Value = 0.8995659351348877
$ws.Cells.Item(41, 4).Value = 0.9240585565567017
$ws.Cells.Item(42, 4).Value = 0.472154289484024
$ws.Cells.Item(43, 4).Value = 0.6990923881530762
$ws.Cells.Item(44, 4).Value = 0.50604248046875
$ws.Cells.Item(45, 4).Value = 0.3106720745563507
$ws.Cells.Item(46, 4).Value = 0.7669169306755066
$ws.Cells.Item(47, 4).Value = 0.7663993835449219
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 0.4958132207393646
$ws.Cells.Item(49, 4).Value = 0.4251040816307068
$ws.Cells.Item(50, 4).Value = 0.3030686676502228
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 0.4973229169845581
$ws.Cells.Item(52, 4).Value = 0.04164242371916771
$ws.Cells.Item(53, 4).Value = 0.4022634327411652
$ws.Cells.Item(54, 4).Value = 0.7426232099533081
$ws.Cells.Item(55, 4).Value = 0.301032692193985
$ws.Cells.Item(56, 4).Value = 0.5794994831085205
$ws.Cells.Item(57, 4).Value = 0.3193793296813965
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 0.4921939074993134
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 0.4344651699066162
$ws.Cells.Item(60, 4).Value = 0.8449845910072327
$ws.Cells.Item(61, 4).Value = 0.660047173500061
$ws.Cells.Item(62, 4).Value = 0.7874123454093933
$ws.Cells.Item(63, 4).Value = 0.601777970790863
$ws.Cells.Item(64, 4).Value = 0.2680604457855225
$ws.Cells.Item(65, 4).Value = 0.1983697563409805
$ws.Cells.Item(66, 4).Value = 0.6706690788269043
$ws.Cells.Item(67, 4).Value = 0.2993455231189728
$ws.Cells.Item(68, 4).Value = 0.6111308932304382
$ws.Cells.Item(69, 4).Value = 0.6785590052604675
$ws.Cells.Item(70, 4).Value = 0.2101347744464874
$ws.Cells.Item(71, 4).Value = 0.3166807889938354
$ws.Cells.Item(72, 4).Value = 0.1608679294586182
$ws.Cells.Item(73, 4).Value = 0.7539152503013611
$ws.Cells.Item(74, 4).Value = 0.4324472844600677
$ws.Cells.Item(75, 4).Value = 0.1621224731206894
$ws.Cells.Item(76, 4).Value = 0.4747454226016998
$ws.Cells.Item(77, 4).Value = 0.5858806371688843
$ws.Cells.Item(78, 4).Value = 0.6373843550682068
$ws.Cells.Item(79, 4).Value = 0.3585901260375977
